$d = $word.ActiveDocument

# 1) Merge the split "Prazo: 28/05/202" + "2" into a single run "Prazo: 28/05/2022"
$d.Content.Find.Execute("Prazo: 28/05/2022", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Prazo: 28/05/2022", 2)

# 2) Fix the typo "Matheus Ferreira" -> "Mateus Ferreira"
$d.Content.Find.Execute("Matheus Ferreira", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Mateus Ferreira", 2)
